# Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column C (td_sim_1); column D mirrors column C for these rows.
$updates = @{
    3  = 92
    6  = 22
    7  = 38
    9  = 18
    11 = 57
    13 = 127
    15 = 243
    17 = 419
    19 = 728
    21 = 14
    23 = 5
    25 = 41
    27 = 65
    29 = 131
    31 = 82
    33 = 30
    35 = 112
    36 = 138
    39 = 2
    40 = 90
    42 = 29
    44 = 20
    46 = 326
    48 = 87
    50 = 45
    52 = 11
    54 = 10
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $ws.Cells.Item($row, 3).Value = $value
    $ws.Cells.Item($row, 4).Value = $value
}

# Recalculate the average in C55 to reflect the corrected td_sim_1 values.
$ws.Cells.Item(55, 3).Value = 110.4444444444444
